# Update countries & provincias Spain
#
# The workbook lists countries (column A) with COVID-19 stats in columns
# B..H ("Casos totales", "Nuevos casos", "Casos activos", "Recuperados",
# "Casos criticos", "Muertes hoy", "Muertes"), one country per row,
# sorted descending by column B ("Casos totales"), starting at row 4
# (rows 1-3 hold the title and header).
#
# This edit refreshes the raw per-country figures for a handful of
# countries and then re-sorts the whole table by "Casos totales"
# (column B) descending, which is what naturally reshuffles the row
# order / shared-string table when the source spreadsheet is
# regenerated. It also bumps the "Datos actualizados ..." timestamp
# cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First row of data and last row of data in the table.
$firstRow = 4
$lastRow = 218

# New raw values (B..H) for every country whose underlying figures
# changed in this update. Countries not listed here keep their current
# values untouched.
$updates = @{
    "Estados Unidos"     = @(1085171, 20977, 150768, 871126, 14781, 1622, 63277)
    "India"              = @(34862,   1800,  9068,   24640,  0,     75,   1154)
    "Israel"             = @(15946,   112,   8561,   7163,   105,   7,    222)
    "Rumania"            = @(12240,   262,   4017,   7506,   221,   24,   717)
    "Guinea"             = @(1495,    144,   329,    1159,   0,     0,    7)
    "Costa Rica"         = @(719,     6,     323,    390,    8,     0,    6)
    "Sudan"              = @(442,     67,    39,     372,    0,     3,    31)
    "Nepal"              = @(57,      0,     16,     8,      0,     0,    0)
    "Republica del Chad" = @(73,      21,    33,     35,     0,     3,    5)
}

foreach ($name in $updates.Keys) {
    $found = $ws.Range("A$firstRow`:A$lastRow").Find($name)
    if ($found -eq $null) {
        Write-Output "WARNING: country not found: $name"
        continue
    }
    $r = $found.Row
    $vals = $updates[$name]
    $ws.Cells.Item($r, 2).Value2 = $vals[0]
    $ws.Cells.Item($r, 3).Value2 = $vals[1]
    $ws.Cells.Item($r, 4).Value2 = $vals[2]
    $ws.Cells.Item($r, 5).Value2 = $vals[3]
    $ws.Cells.Item($r, 6).Value2 = $vals[4]
    $ws.Cells.Item($r, 7).Value2 = $vals[5]
    $ws.Cells.Item($r, 8).Value2 = $vals[6]
}

# Re-sort the whole data table by "Casos totales" (column B), descending,
# so the row order (and therefore the shared-string table order) matches
# the refreshed figures.
$sortRange = $ws.Range("A$firstRow`:H$lastRow")
$sortKey = $ws.Range("B$firstRow`:B$lastRow")
$sortRange.Sort($sortKey, 2)

# Bump the "last updated" timestamp.
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 21:52"

Write-Output "Update complete"
